$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translation")

$rows = @(
    ,@(736, 'SingleUseId835', 'Typography_05', 'Left', 'LTR', 'New Text')
    ,@(737, 'SingleUseId836', 'Typography_05', 'Left', 'LTR', 'ECU information')
    ,@(738, 'SingleUseId837', 'Typography_01', 'Left', 'LTR', 'ACTIVATING')
    ,@(739, 'SingleUseId838', 'Typography_01', 'Left', 'LTR', 'OK')
    ,@(740, 'SingleUseId839', 'Typography_01', 'Left', 'LTR', 'CHARGING')
    ,@(741, 'SingleUseId840', 'Typography_01', 'Left', 'LTR', 'UNKNOWN')
    ,@(742, 'SingleUseId841', 'Typography_01', 'Left', 'LTR', 'FATAL')
    ,@(743, 'SingleUseId842', 'Typography_01', 'Left', 'LTR', 'BMS HV')
    ,@(744, 'SingleUseId843', 'Typography_01', 'Left', 'LTR', 'BMS LV')
    ,@(745, 'SingleUseId844', 'Typography_01', 'Left', 'LTR', 'APPS')
    ,@(746, 'SingleUseId845', 'Typography_01', 'Left', 'LTR', 'INVERTER')
    ,@(747, 'SingleUseId846', 'Typography_01', 'Left', 'LTR', 'DOF')
    ,@(748, 'SingleUseId847', 'Typography_01', 'Left', 'LTR', 'GPS')
    ,@(749, 'SingleUseId848', 'Typography_01', 'Left', 'LTR', 'SD')
    ,@(750, 'SingleUseId849', 'Typography_01', 'Left', 'LTR', 'XBEE')
    ,@(751, 'SingleUseId850', 'Typography_02', 'Left', 'LTR', 'CELL 5')
    ,@(752, 'SingleUseId851', 'Typography_02', 'Left', 'LTR', 'CELL 6')
    ,@(753, 'SingleUseId856', 'Typography_05', 'Center', 'LTR', '<value>C')
    ,@(754, 'SingleUseId857', 'Typography_05', 'Left', 'LTR', '''0')
    ,@(755, 'SingleUseId858', 'Typography_05', 'Center', 'LTR', '<value>C')
    ,@(756, 'SingleUseId859', 'Typography_05', 'Left', 'LTR', '''0')
    ,@(757, 'SingleUseId868', 'Typography_05', 'Left', 'LTR', 'Precharg info')
    ,@(758, 'SingleUseId869', 'Typography_05', 'Left', 'LTR', 'finished')
    ,@(759, 'SingleUseId870', 'Typography_05', 'Left', 'LTR', 'nije se zatvorio A-')
    ,@(760, 'SingleUseId872', 'Default', 'Left', 'LTR', 'AIR-')
    ,@(761, 'SingleUseId873', 'Default', 'Left', 'LTR', 'AIR+')
    ,@(762, 'SingleUseId874', 'Typography_05', 'Left', 'LTR', 'Shutdown info')
    ,@(763, 'SingleUseId875', 'Typography_05', 'Left', 'LTR', 'open')
    ,@(764, 'SingleUseId876', 'Default', 'Left', 'LTR', 'Relay')
    ,@(765, 'SingleUseId877', 'Default', 'Left', 'LTR', 'SHD cmd')
    ,@(766, 'SingleUseId878', 'Default', 'Left', 'LTR', 'IMD')
    ,@(767, 'SingleUseId879', 'Default', 'Left', 'LTR', 'Current')
    ,@(768, 'SingleUseId880', 'Typography_05', 'Left', 'LTR', 'nije se zatvorio A+')
    ,@(769, 'SingleUseId881', 'Default', 'Left', 'LTR', 'Voltage')
    ,@(770, 'SingleUseId882', 'Typography_05', 'Left', 'LTR', 'not finished')
    ,@(771, 'SingleUseId883', 'Typography_05', 'Left', 'LTR', 'closed')
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
}

